# Apply workbook-level and sheet-level changes described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Sheet1" -> "PD_ref_cat"
$ws.Name = "PD_ref_cat"

# Make sure the renamed sheet is the active sheet/tab, then move the
# current selection from A4 to A9 (as reflected in the sheetView's
# <selection activeCell="A9" sqref="A9"/>).
$ws.Activate()
$ws.Range("A9").Select()
